$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $val) {
    $origStyle = $cellRange.Style
    $cellRange.Value = "'" + $val
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "66.449.58"
$ws.Range("E2").Value = "  +0.43%  "
Set-TextValue $ws.Range("D3") "3.251.64"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "610.15"
$ws.Range("E5").Value = "  +0.82%  "
Set-TextValue $ws.Range("D6") "157.31"
$ws.Range("E6").Value = "  +2.23%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.06%  "
Set-TextValue $ws.Range("D8") "3.252.17"
$ws.Range("E8").Value = "  +2.74%  "
$ws.Range("E9").Value = "  -0.12%  "
Set-TextValue $ws.Range("D10") "0.162"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  +4.01%  "
Set-TextValue $ws.Range("D12") "0.499"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("E13").Value = "  +1.30%  "
Set-TextValue $ws.Range("D14") "39.06"
$ws.Range("E14").Value = "  +1.94%  "
Set-TextValue $ws.Range("D15") "3.780.70"
$ws.Range("E15").Value = "  +2.55%  "
Set-TextValue $ws.Range("D16") "66.535.12"
$ws.Range("E16").Value = "  +0.51%  "
Set-TextValue $ws.Range("D17") "7.48"
$ws.Range("E17").Value = "  +0.85%  "
Set-TextValue $ws.Range("D18") "3.246.09"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("E19").Value = "  +1.24%  "
Set-TextValue $ws.Range("D20") "504.72"
$ws.Range("E20").Value = "  -0.99%  "
Set-TextValue $ws.Range("D21") "15.45"
$ws.Range("E21").Value = "  +0.44%  "
Set-TextValue $ws.Range("D22") "0.752"
$ws.Range("E22").Value = "  +3.40%  "
Set-TextValue $ws.Range("D23") "8.10"
$ws.Range("E23").Value = "  +1.15%  "
Set-TextValue $ws.Range("D24") "14.72"
$ws.Range("E24").Value = "  +0.00%  "
Set-TextValue $ws.Range("D25") "87.18"
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.66%  "
Set-TextValue $ws.Range("D28") "9.17"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E29").Value = "  +0.38%  "
Set-TextValue $ws.Range("D30") "0.128"
$ws.Range("E30").Value = "  +45.11%  "
$ws.Range("E31").Value = "  -2.36%  "
Set-TextValue $ws.Range("D32") "2.89"
$ws.Range("E32").Value = "  -4.89%  "
Set-TextValue $ws.Range("D33") "27.95"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  -3.46%  "
Set-TextValue $ws.Range("D36") "6.45"
$ws.Range("E36").Value = "  -0.74%  "
Set-TextValue $ws.Range("D37") "55.51"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").Value = "  +18.77%  "
Set-TextValue $ws.Range("D39") "0.0₃0784"
$ws.Range("E39").Value = "  +15.16%  "
Set-TextValue $ws.Range("D40") "495.82"
$ws.Range("E40").Value = "  -0.94%  "
Set-TextValue $ws.Range("D41") "0.0423"
$ws.Range("E41").Value = "  +0.93%  "
Set-TextValue $ws.Range("D42") "0.129"
$ws.Range("E42").Value = "  +0.63%  "
Set-TextValue $ws.Range("D43") "8.84"
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D44") "2.52"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D45") "0.293"
$ws.Range("E45").Value = "  -0.76%  "
Set-TextValue $ws.Range("D46") "2.988.20"
$ws.Range("E46").Value = "  +5.80%  "
Set-TextValue $ws.Range("D47") "28.82"
$ws.Range("E47").Value = "  +3.17%  "
Set-TextValue $ws.Range("D48") "2.53"
$ws.Range("E48").Value = "  +6.48%  "
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
Set-TextValue $ws.Range("D51") "2.53"
$ws.Range("E51").Value = "  -2.97%  "
